$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (cascades automatically to the _FilterDatabase defined name)
$ws.Name = "DoesNotMatter"

# Move the active selection to M27
$ws.Activate()
$ws.Range("M27").Select()
